$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added remove item functionality" -> the visible data change is: 10 more
# order rows got added (rows 21-30), and the previous "lmao" comment string
# is repurposed/renamed to "20" for those new rows, while the four existing
# rows that used to say "lmao" (rows 17-20) are reset to a blank " " comment.

# 1) Duplicate the "FRIES / NTU / side" block (rows 2-6) twice to create the
#    10 new rows 21-30, then renumber their Order ID to 4.
$ws.Range("A2:F6").Copy($ws.Range("A21:F25"))
$ws.Range("A2:F6").Copy($ws.Range("A26:F30"))
$ws.Range("A21:A30").Value = 4

# 2) The new rows' Comments column (F) should hold the TEXT "20" (not the
#    number 20), so build it via a TEXT() formula in a scratch cell, copy
#    the resulting text, and paste the value into each new row's F cell.
$ws.Range("H1").Formula = '=TEXT(20,"0")'
$ws.Range("H1").Copy()
for ($r = 21; $r -le 30; $r++) {
    $ws.Cells.Item($r, 6).PasteSpecial(-4163)
}
$ws.Range("H1").ClearContents()

# 3) The old rows (17-20) that used to say "lmao" now just have a blank " "
#    comment.
$ws.Range("F17:F20").Replace("lmao", " ")
